$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 new columns before column B ---
# Old layout: B=Jun_17, C=Jun_15, D=Jun_13, E=Jun_10
# New layout: B=Jun_27, C=Jun_26, D=Jun_26, E=Jun_17, F=Jun_15, G=Jun_13, H=Jun_10
$ws.Range("B1:D1").EntireColumn.Insert() | Out-Null

# --- Header row (row 1) for the newly inserted columns ---
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# --- Fill data rows (2-27) for new columns B:D with "UN" ---
$ws.Range("B2:D27").Value = "UN"

# --- Add two new rows at the bottom for the new analyst group ---
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28:D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29:D29").Value = "UN"

Write-Output "Edit complete"
